# Finished all website scraping: append new website names to column A
# for the rows that already hold printer model data in column B (rows 14-17),
# and add a brand new row (18) for the last scraped website.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A already contains: Staples, NewEgg, B&H, Walmart, BuyVC, Dell,
# Vistek, Tastar, Overland, PC Nation, HP, Tiger Direct (rows 2-13).
# Continue the list with websites that were already scraped into the
# shared-string pool for the "Blocked Websites" column (Adorama, GoVets,
# Plotter), followed by two brand-new website names.
$ws.Range("A14").Value = "Adorama"
$ws.Range("A15").Value = "GoVets"
$ws.Range("A16").Value = "Plotter"
$ws.Range("A17").Value = "PC Connection"
$ws.Range("A18").Value = "Amazon"

# Move the active selection to reflect where editing finished.
$ws.Range("B18").Select()
